# Binary Search: Square root of N upto 3 decimal places - updated index
#
# The workbook has two sheets; "Binary Search 2" (index 2) holds the index
# table for the binary-search folder. A new row entry (row 4, question #2)
# already has its first few columns filled in; this edit adds the GitHub
# link for the solution file in column F, grows the row to fit the wrapped
# link text, and moves the active selection down to the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

# Add the GitHub hyperlink to F4: sets the cell value/display text and
# registers the external relationship.
$ws.Hyperlinks.Add(
    $ws.Range("F4"),
    "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/2_square_root_of_n_upto_3_decimal_places.java",
    [Type]::Missing,
    [Type]::Missing,
    "dsa/2_square_root_of_n_upto_3_decimal_places.java at main · ankurnecessary/dsa · GitHub"
) | Out-Null

# Give F4 the same look as the existing hyperlink cell F3 (wrapped-text
# Hyperlink style) by copying its formatting over.
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The extra wrapped line in F4 makes row 4 taller.
$ws.Rows.Item(4).RowHeight = 57.6

# Move the active selection to the next row, ready for further entries.
$ws.Range("B5").Select() | Out-Null
